$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Internal control"
$ws.Range("A10").Value = "Not clear instruction"
$ws.Range("A11").Value = "Drawing update"

$ws.Range("A12").Select()
